$d = $word.ActiveDocument

# 1 & 2: Update "Trend difference" table header text from 2014-01 to 2015-01
# (both occurrences in the table get replaced by a single wrapped Find/Replace)
$d.Content.Find.Execute(
    "Trend difference (2014-01 to 2022-12)", $true, $false, $false, $false,
    $false, $true, 1, $false, "Trend difference (2015-01 to 2022-12)", 2)

# 3: Remove the unused "Abstract Title" paragraph style
$abstractTitle = $d.Styles("AbstractTitle")
$abstractTitle.Delete()

# 4: Change the "Abstract" style's space-before from 100 twips (5pt) to 300 twips (15pt)
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 15

# 5: Remove the unused "Footnote Block Text" paragraph style
$footnoteBlockText = $d.Styles("FootnoteBlockText")
$footnoteBlockText.Delete()
